$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.743.34'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '2.277.65'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.67'
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('E6').Value = '  +2.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '75.14'
$ws.Range('E7').Value = '  +6.77%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.635'
$ws.Range('E9').Value = '  -4.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.34'
$ws.Range('E10').Value = '  +0.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0972'
$ws.Range('E11').Value = '  +0.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.41'
$ws.Range('E12').Value = '  -0.74%  '
$ws.Range('E13').Value = '  +0.65%  '
$ws.Range('D14').Value = '2.618.55'
$ws.Range('E14').Value = '  +1.54%  '
$ws.Range('E15').Value = '  +1.73%  '
$ws.Range('E16').Value = '  -1.49%  '
$ws.Range('D17').Value = '2.274.59'
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('D18').Value = '42.638.78'
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').Value = '0.0₃0997'
$ws.Range('E19').Value = '  +0.91%  '
$ws.Range('E20').Value = '  -1.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.42'
$ws.Range('E21').Value = '  -0.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '234.66'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('E23').Value = '  +6.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.87'
$ws.Range('E24').Value = '  -1.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.25'
$ws.Range('E26').Value = '  -1.58%  '
$ws.Range('E27').Value = '  -1.13%  '
$ws.Range('E28').Value = '  +1.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '167.49'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '21.04'
$ws.Range('E30').Value = '  +1.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0872'
$ws.Range('E31').Value = '  +9.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.34'
$ws.Range('E32').Value = '  -1.53%  '
$ws.Range('E33').Value = '  +0.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '31.84'
$ws.Range('E34').Value = '  +1.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.127'
$ws.Range('E35').Value = '  +1.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.55'
$ws.Range('E36').Value = '  +2.33%  '
$ws.Range('E37').Value = '  +1.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0304'
$ws.Range('E38').Value = '  -5.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.55'
$ws.Range('E39').Value = '  +8.84%  '
$ws.Range('E40').Value = '  -0.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.87'
$ws.Range('E41').Value = '  +1.44%  '
$ws.Range('E42').Value = '  +3.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '61.29'
$ws.Range('E43').Value = '  -1.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.91'
$ws.Range('E44').Value = '  -0.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '105.76'
$ws.Range('E45').Value = '  +11.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.75'
$ws.Range('E46').Value = '  -2.90%  '
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('E49').Value = '  -0.49%  '
$ws.Range('E50').Value = '  -1.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.21'
$ws.Range('E51').Value = '  -0.68%  '
